$d = $word.ActiveDocument

# The "Accounts" section has two tables; the second (word/definition legend
# table is first) holds the actual account rows we need to touch.
$tbl = $d.Tables.Item(2)

# 1) gocomics.com -> https://gocomics.com  (row 3, "Account" column)
$tbl.Cell(3, 1).Range.Text = "https://gocomics.com"

# 2) comicskingdom.com -> https://comicskingdom.com  (row 4, "Account" column)
$tbl.Cell(4, 1).Range.Text = "https://comicskingdom.com"

# 3) comicskingdom.com row: Subscription column (col 4) flips No -> Yes
$tbl.Cell(4, 4).Range.Text = "Yes"

# 4) Append a brand-new row for the stanza.co account
$tbl.Rows.Add() | Out-Null
$newRowIndex = $tbl.Rows.Count
$tbl.Cell($newRowIndex, 1).Range.Text = "https://stanza.co"
$tbl.Cell($newRowIndex, 2).Range.Text = "ralph.hightower@gmail.com"
$tbl.Cell($newRowIndex, 3).Range.Text = "No"
$tbl.Cell($newRowIndex, 4).Range.Text = "Yes"
$tbl.Cell($newRowIndex, 5).Range.Text = "Cancel"

Write-Output "Accounts table updated"
